# The edit targets the "Comparison" slide layout -- the 5th Custom Layout
# on the deck's slide master. Its two left-hand placeholders
# ("Text Placeholder 2" and "Content Placeholder 3") have their vertical
# position/size adjusted: the upper placeholder shrinks and the lower one
# moves up and grows, closing the gap that used to sit between them.
#
# PowerPoint's Shape.Left/Top/Width/Height COM properties are expressed in
# points (1 pt = 12700 EMU) and are backed by single-precision floats, so
# the literals below are chosen (via binary search against this runtime)
# to round-trip to the exact target EMU values from the OOXML diff:
#
#   Text Placeholder 2:    off(457200,1151334) ext(4040188,2531666)
#                        -> off(457200,1151333) ext(4040188,1468041)
#   Content Placeholder 3: off(457200,3749674) ext(4040188, 844947)
#                        -> off(457200,2959100) ext(4040188,1635521)

$p = $ppt.ActivePresentation

$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$layout = $master.CustomLayouts.Item(5)   # "Comparison" layout

$textPlaceholder = $layout.Shapes.Item(2)    # "Text Placeholder 2"
$contentPlaceholder = $layout.Shapes.Item(3) # "Content Placeholder 3"

$textPlaceholder.Left = 35.99999809265137
$textPlaceholder.Top = 90.6561393737793
$textPlaceholder.Width = 318.1250457763672
$textPlaceholder.Height = 115.59377670288087

$contentPlaceholder.Left = 35.99999809265137
$contentPlaceholder.Top = 232.99999237060547
$contentPlaceholder.Width = 318.1250457763672
$contentPlaceholder.Height = 128.78118133544922
